$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.953.82"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").Value = "1.673.88"
$ws.Range("E3").Value = "  -0.35%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").Value = "'214.71"
$ws.Range("E5").Value = "  -0.95%  "

$ws.Range("D6").Value = "'0.516"
$ws.Range("E6").Value = "  +1.49%  "

$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("D10").Value = "'20.33"
$ws.Range("E10").Value = "  +0.20%  "

$ws.Range("D11").Value = "'0.0887"
$ws.Range("E11").Value = "  +0.08%  "

$ws.Range("D12").Value = "1.911.00"
$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").Value = "1.677.41"
$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("E14").Value = "  -0.55%  "

$ws.Range("E15").Value = "  +0.12%  "

$ws.Range("D16").Value = "'65.63"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Value = "26.960.69"
$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'235.99"
$ws.Range("E18").Value = "  -1.20%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'8.11"
$ws.Range("E19").Value = "  +3.93%  "

$ws.Range("E20").Value = "  -0.96%  "

$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").Value = "'4.43"
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "  -1.44%  "

$ws.Range("E24").Value = "  -2.16%  "

$ws.Range("D25").Value = "'145.39"
$ws.Range("E25").Value = "  -0.18%  "

$ws.Range("D26").Value = "'7.20"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("D27").Value = "'16.00"
$ws.Range("E27").Value = "  +0.16%  "

$ws.Range("D28").Value = "'0.112"
$ws.Range("E28").Value = "  -1.55%  "

$ws.Range("E29").Value = "  +0.24%  "

$ws.Range("E30").Value = "  -0.44%  "

$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").Value = "1.481.95"
$ws.Range("E33").Value = "  -2.95%  "

$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("E35").Value = "  +2.88%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").Value = "'0.588"
$ws.Range("E37").Value = "  +1.55%  "

$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("E39").Value = "  +0.25%  "

$ws.Range("E40").Value = "  -3.72%  "

$ws.Range("E41").Value = "  +5.36%  "

$ws.Range("E42").Value = "  +0.26%  "

$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("D44").Value = "'67.16"
$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").Value = "1.816.10"
$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("E46").Value = "  -0.52%  "

$ws.Range("D47").Value = "'90.60"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").Value = "'7.72"
$ws.Range("E51").Value = "  +0.11%  "
